$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.554.75"
$ws.Range("E2").Value = "  +2.04%  "
# Row 3
$ws.Range("D3").Value = "1.893.46"
$ws.Range("E3").Value = "  +1.72%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.025"
$ws.Range("E4").Value = "  +2.14%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.52"
$ws.Range("E5").Value = "  +2.07%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.024"
$ws.Range("E6").Value = "  +2.17%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5168"
$ws.Range("E7").Value = "  +0.61%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3929"
$ws.Range("E8").Value = "  +2.68%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08296"
$ws.Range("E9").Value = "  +0.31%  "
# Row 10
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.127"
$ws.Range("E10").Value = "  +1.66%  "
# Row 11
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.921.57"
$ws.Range("E11").Value = "  +3.85%  "
# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.260"
$ws.Range("E12").Value = "  +1.09%  "
# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.36"
$ws.Range("E13").Value = "  -0.69%  "
# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.260"
$ws.Range("E14").Value = "  -0.13%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.023"
$ws.Range("E15").Value = "  +1.96%  "
# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001107"
$ws.Range("E16").Value = "  +1.10%  "
# Row 17
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.24"
$ws.Range("E17").Value = "  +0.88%  "
# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06781"
$ws.Range("E18").Value = "  +2.07%  "
# Row 19
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.84"
$ws.Range("E19").Value = "  +1.00%  "
# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.024"
$ws.Range("E20").Value = "  +2.12%  "
# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.044"
$ws.Range("E21").Value = "  +0.45%  "
# Row 22
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "28.605.92"
$ws.Range("E22").Value = "  +2.11%  "
# Row 23
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.16"
$ws.Range("E23").Value = "  +0.83%  "
# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.263"
$ws.Range("E24").Value = "  +0.87%  "
# Row 25
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.413"
$ws.Range("E25").Value = "  +0.79%  "
# Row 26
$ws.Range("D26").Value = "2.118.42"
$ws.Range("E26").Value = "  +2.43%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.54"
$ws.Range("E27").Value = "  +2.61%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.87"
$ws.Range("E28").Value = "  +1.89%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.440"
$ws.Range("E29").Value = "  -2.43%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.92"
$ws.Range("E30").Value = "  +1.86%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1061"
$ws.Range("E31").Value = "  -0.11%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.044"
$ws.Range("E32").Value = "  +1.70%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.920"
$ws.Range("E33").Value = "  +1.86%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.666"
$ws.Range("E34").Value = "  +2.05%  "
# Row 35
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.362"
$ws.Range("E35").Value = "  -0.83%  "
# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02447"
$ws.Range("E36").Value = "  +1.72%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06597"
$ws.Range("E37").Value = "  +1.68%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2195"
$ws.Range("E38").Value = "  +0.27%  "
# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.267"
$ws.Range("E39").Value = "  +4.90%  "
# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6516"
$ws.Range("E40").Value = "  -0.29%  "
# Row 41
$ws.Range("E41").Value = "  -0.87%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.983"
$ws.Range("E42").Value = "  +0.10%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.17"
$ws.Range("E43").Value = "  -0.07%  "
# Row 44
$ws.Range("E44").Value = "  +0.13%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.16"
$ws.Range("E45").Value = "  +1.13%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.740"
$ws.Range("E46").Value = "  +1.88%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.296"
$ws.Range("E47").Value = "  +1.31%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.011"
$ws.Range("E48").Value = "  -0.20%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.223"
$ws.Range("E49").Value = "  +0.99%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.80"
$ws.Range("E50").Value = "  +1.74%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06941"
$ws.Range("E51").Value = "  +1.94%  "
